# The document contains four paragraphs with the (stale, 2018 Perseus)
# observation-window sentence, spread across several runs with explicit
# character formatting (rFonts/lang). Per the commit, each of these
# paragraphs is collapsed down to a single, plain run (no rPr at all)
# containing the new, not-yet-fully-translated sentence about Hercules.
#
# Strategy: for each target paragraph, delete the run content (but not
# the paragraph mark, so the paragraph/pPr survives), then InsertAfter()
# on the now-empty paragraph range — InsertAfter on an empty range creates
# a brand-new run with no inherited rPr, matching the diff exactly.

$d = $word.ActiveDocument

$oldText = "V roku 2018 môžete pozorovať súhvezdie Perseus: 30. októbra - 8. novembra a 29. novembra - 8. decembra"
$newText = "V roku Hercules: 13.-22. júna, 12.-21. júla, 10.-19. augusta"

$targets = New-Object System.Collections.ArrayList

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    $trimmed = $t.Trim()
    if ($trimmed -eq $oldText) {
        [void]$targets.Add($i)
    }
}

foreach ($idx in $targets) {
    $full = $d.Paragraphs($idx).Range
    # Exclude the trailing paragraph mark so the <w:p>/pPr survives.
    $body = $d.Range($full.Start, $full.End - 1)
    $body.Delete()

    $p = $d.Paragraphs($idx).Range
    $p.InsertAfter($newText)
}
